$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the repayment strategy value from "Mifos style" to the new scenario text,
# matching the formatting already used by B1 (left/top aligned, same fill/font).
$ws.Range("B1").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B17").Value = "Penalties, Fees, Interest, Principal order"

# Move the active selection to B17 to match the saved cursor position.
$ws.Range("B17").Select()
